$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "94.403.25"
Set-TextValue "E2" "  -3.90%  "
Set-TextValue "D3" "3.428.59"
Set-TextValue "E3" "  +1.21%  "
Set-TextValue "E4" "  +0.04%  "
Set-TextValue "D5" "237.82"
Set-TextValue "E5" "  -6.24%  "
Set-TextValue "D6" "642.47"
Set-TextValue "E6" "  -2.96%  "
Set-TextValue "D7" "1.45"
Set-TextValue "E7" "  -0.60%  "
Set-TextValue "D8" "0.406"
Set-TextValue "E8" "  -4.43%  "
Set-TextValue "E9" "  +0.10%  "
Set-TextValue "D10" "0.979"
Set-TextValue "E10" "  -6.47%  "
Set-TextValue "D11" "3.427.82"
Set-TextValue "E11" "  +1.26%  "
Set-TextValue "E12" "  -4.31%  "
Set-TextValue "D13" "42.02"
Set-TextValue "E13" "  +0.73%  "
Set-TextValue "D14" "6.24"
Set-TextValue "E14" "  +1.83%  "
Set-TextValue "D15" "94.178.67"
Set-TextValue "E15" "  -3.76%  "
Set-TextValue "D16" "4.066.90"
Set-TextValue "E16" "  +1.26%  "
Set-TextValue "D17" "0.0000252"
Set-TextValue "E17" "  -1.47%  "
Set-TextValue "D18" "8.39"
Set-TextValue "E18" "  -5.88%  "
Set-TextValue "D19" "3.425.24"
Set-TextValue "E19" "  +1.22%  "
Set-TextValue "D20" "17.56"
Set-TextValue "E20" "  -2.70%  "
Set-TextValue "D21" "11.58"
Set-TextValue "E21" "  +5.66%  "
Set-TextValue "D22" "0.497"
Set-TextValue "E22" "  -6.00%  "
Set-TextValue "D23" "500.15"
Set-TextValue "E23" "  -2.49%  "
Set-TextValue "D24" "3.24"
Set-TextValue "E24" "  -6.10%  "
Set-TextValue "E25" "  -4.42%  "
Set-TextValue "D26" "6.54"
Set-TextValue "E26" "  -6.36%  "
Set-TextValue "D27" "94.25"
Set-TextValue "E27" "  -2.53%  "
Set-TextValue "D28" "12.00"
Set-TextValue "E28" "  -2.91%  "
Set-TextValue "D29" "3.609.58"
Set-TextValue "E29" "  +1.16%  "
Set-TextValue "D30" "11.79"
Set-TextValue "E30" "  +3.25%  "
Set-TextValue "E31" "  -0.09%  "
Set-TextValue "E32" "  +6.51%  "
Set-TextValue "D33" "0.139"
Set-TextValue "E33" "  -2.94%  "
Set-TextValue "D34" "1.00"
Set-TextValue "E34" "  -0.08%  "
Set-TextValue "D35" "0.180"
Set-TextValue "E35" "  -4.20%  "
Set-TextValue "D36" "29.64"
Set-TextValue "E36" "  +2.15%  "
Set-TextValue "D37" "0.553"
Set-TextValue "E37" "  -1.33%  "
Set-TextValue "D38" "560.61"
Set-TextValue "E38" "  +4.33%  "
Set-TextValue "D39" "1.46"
Set-TextValue "E39" "  -3.25%  "
Set-TextValue "D40" "7.65"
Set-TextValue "E40" "  -4.57%  "
Set-TextValue "E41" "  +0.00%  "
Set-TextValue "D42" "0.151"
Set-TextValue "E42" "  -1.19%  "
Set-TextValue "D43" "0.907"
Set-TextValue "E43" "  +5.87%  "
Set-TextValue "D44" "24.06"
Set-TextValue "E44" "  -1.47%  "
Set-TextValue "D45" "1.73"
Set-TextValue "E45" "  +0.19%  "
Set-TextValue "B46" "Filecoin"
Set-TextValue "C46" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D46" "5.68"
Set-TextValue "E46" "  +1.24%  "
Set-TextValue "B47" "MantraDAO"
Set-TextValue "C47" "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
Set-TextValue "D47" "3.65"
Set-TextValue "E47" "  -0.76%  "
Set-TextValue "B48" "VeChain"
Set-TextValue "C48" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D48" "0.0412"
Set-TextValue "E48" "  -3.95%  "
Set-TextValue "B49" "OKB"
Set-TextValue "C49" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D49" "55.24"
Set-TextValue "E49" "  -1.57%  "
Set-TextValue "D50" "3.34"
Set-TextValue "E50" "  +3.49%  "
Set-TextValue "D51" "2.19"
Set-TextValue "E51" "  -3.62%  "
